$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report date (keep it stored as text, matching the original cell type)
$ws.Range("D1").NumberFormat = "@"
$ws.Range("D1").Value = "2022-09-21"

# Rename "Ukupan broj kafa: " label to "Ukupan broj torti: "
$ws.Range("A7").Value = "Ukupan broj torti: "

# Update the most-commented post title
$ws.Range("B15").Value = "Čoko-lešnik torta (Novogodišnja torta)"

# Update numeric statistics
$ws.Range("B3").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("B7").Value = 2
$ws.Range("B9").Value = 4
$ws.Range("B11").Value = 1
$ws.Range("B13").Value = 0.5
